# KostMan docs/clean.xlsx - room details / meter_start column fix
#
# The commit adds a new "meter_start" data column (Q) that was previously
# only present as a header (Q1) but had no values underneath it. Populate
# Q2:Q174 with the initial/default meter-start reading of 0 for every room
# record row, matching the header already present in Q1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sanity: make sure the header is what we expect before filling the column.
# (Q1 already contains "meter_start" in the source workbook.)

$lastRow = 174
$qRange = "Q2:Q" + $lastRow

# Fill the new meter_start values for every data row (2-174) with 0.
$ws.Range($qRange).Value = 0

# Reflect the editor's final on-screen selection/viewport for this sheet:
# the whole new Q column (Q2:Q174) ends up selected, with the view scrolled
# so column M is the left-most visible column and row 1 is visible.
$ws.Activate()
$ws.Range("M1").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Range($qRange).Select()
